$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text.TrimEnd("`r", "`a")
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1. After "Takes in arguments and performs UML parsing" paragraph,
#    insert a new sub-bullet: "Pass in argument and search for path"
# ---------------------------------------------------------------------
$pTakesIn = Find-ParagraphByText $d "Takes in arguments and performs UML parsing"
$rng = $pTakesIn.Range
$rng.Collapse(0)
$rng.InsertParagraphAfter()
$pNew1 = $pTakesIn.Next()
$pNew1.Range.Text = "Pass in argument and search for path"
$pNew1.Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------
# 2. Merge the "Create an output file..." paragraph's three runs
#    (incl. the spell-checked "Umple") into a single run, dropping the
#    proofErr spell-check markers.
# ---------------------------------------------------------------------
$pCreate = Find-ParagraphByText $d "Create an output file (in text) to store the Umple code"
$rngCreate = $pCreate.Range
$rngCreate.MoveEnd(1, -1)
$rngCreate.Delete()
$rngCreate.InsertAfter("Create an output file (in text) to store the Umple code")

# ---------------------------------------------------------------------
# 3. After "Create an output file..." paragraph, insert a new
#    sub-bullet: "Pass in path file to output file. ..."
# ---------------------------------------------------------------------
$pCreate2 = Find-ParagraphByText $d "Create an output file (in text) to store the Umple code"
$rngCreate2 = $pCreate2.Range
$rngCreate2.Collapse(0)
$rngCreate2.InsertParagraphAfter()
$pNew2 = $pCreate2.Next()
$pNew2.Range.Text = "Pass in path file to output file.  Default behavior should output it in the same folder as the file to be parsed"
$pNew2.Range.ListFormat.ListLevelNumber = 2

# ---------------------------------------------------------------------
# 4. Split "Checks the methods and what the type of methods are
#    (void or return type)" into five separate runs with new wording.
# ---------------------------------------------------------------------
$pChecks = Find-ParagraphByText $d "Checks the methods and what the type of methods are (void or return type)"
$rngChecks = $pChecks.Range
$rngChecks.MoveEnd(1, -1)
$rngChecks.Delete()
$rngChecks2 = $pChecks.Range
$rngChecks2.MoveEnd(1, -1)
$xmlFrag = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Checks the methods and </w:t></w:r><w:r><w:t xml:space="preserve">their return </w:t></w:r><w:r><w:t>type</w:t></w:r><w:r><w:t xml:space="preserve">s </w:t></w:r><w:r><w:t>(void or return type)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rngChecks2.InsertXML($xmlFrag)

# ---------------------------------------------------------------------
# 5. "Once finished with reading" -> "Once finished reading, output the
#    code into a diagram (optional)", and promote from ilvl=1 to ilvl=0
# ---------------------------------------------------------------------
$pOnce = Find-ParagraphByText $d "Once finished with reading"
$rngOnce = $pOnce.Range
$rngOnce.ListFormat.ListLevelNumber = 1
$find = $rngOnce.Find
$find.Execute("Once finished with reading", $true, $false, $false, $false, $false, $true, 1, $false, "Once finished reading, output the code into a diagram (optional)", 2)
